# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-row labels to "_FV2410" / "_FV2504"
# - Freeze the header row (row 1)
# - Turn the used range into an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels --------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# --- 2. Freeze the top row ---------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the used range into a Table -----------------------------
$rng = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
